# Trade #108 closed at 2026-02-17 09:19:03 - unknown UNKNOWN +0.000%
#
# Updates the summary stats to reflect the newly closed trade and appends
# the trade row itself to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - roll the aggregate stats forward
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.1    # Current Capital
$summary.Range("B4").Value = 0.11      # Total P&L $
$summary.Range("B5").Value = 0.02      # Total P&L %
$summary.Range("B6").Value = 108       # Total Trades
$summary.Range("B8").Value = 43        # Losing Trades
$summary.Range("B9").Value = 42.59     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - same roll-forward for the MarketMaking row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.1      # Capital
$status.Range("D4").Value = 108        # Trades
$status.Range("E4").Value = 0.11       # P&L $
$status.Range("F4").Value = 0.1        # P&L %
$status.Range("G4").Value = 42.59      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade row (#108 / row 109) to "All Trades" and
#    "MarketMaking". Date/Time columns must stay plain text, so force a
#    text number format before writing them to avoid Excel's automatic
#    date/time conversion.
# ---------------------------------------------------------------------
$tradeRow = 109

function Add-TradeRow($ws) {
    # Force text format on the Date/Time cells first so Excel doesn't
    # auto-convert the literal strings into date/time serial numbers, then
    # restore the plain "Normal" style so no stray number format lingers.
    $dateTimeRange = $ws.Range("B" + $tradeRow + ":C" + $tradeRow)
    $dateTimeRange.NumberFormat = "@"

    $ws.Cells.Item($tradeRow, 1).Value = 108                                         # Trade #
    $ws.Cells.Item($tradeRow, 2).Value = "2026-02-17"                                # Date
    $ws.Cells.Item($tradeRow, 3).Value = "09:18:57"                                  # Time

    $dateTimeRange.Style = "Normal"

    $ws.Cells.Item($tradeRow, 4).Value = "MarketMaking"                              # Strategy
    $ws.Cells.Item($tradeRow, 5).Value = "DOWN"                                      # Side
    $ws.Cells.Item($tradeRow, 6).Value = 0.83                                        # Entry Price
    $ws.Cells.Item($tradeRow, 7).Value = 0.8                                         # Exit Price
    $ws.Cells.Item($tradeRow, 8).Value = "CLOSED"                                    # Status
    $ws.Cells.Item($tradeRow, 9).Value = -3.6145                                     # P&L %
    $ws.Cells.Item($tradeRow, 10).Value = -0.03                                      # P&L $
    $ws.Cells.Item($tradeRow, 11).Value = 100.1                                      # Capital After
    $ws.Cells.Item($tradeRow, 12).Value = 0                                          # Entry Slippage (bps)
    $ws.Cells.Item($tradeRow, 13).Value = 0                                          # Exit Slippage (bps)
    $ws.Cells.Item($tradeRow, 14).Value = 0.6                                        # Confidence
    $ws.Cells.Item($tradeRow, 15).Value = "Normal spread capture: 19600 bps"         # Entry Reason
    $ws.Cells.Item($tradeRow, 16).Value = "early_exit"                               # Exit Reason
    $ws.Cells.Item($tradeRow, 17).Value = 0.14                                       # Duration (min)
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
